$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: naive forecaster bug fix drops the trailing (10th) error column; no new leading value for this row
$ws.Range("K2").ClearContents()

# Row 3: insert new first-quarter error in B3; shift prior values right by one column
$ws.Range("B3").Value = [double]"1.052939957446597E-10"
$ws.Range("C3").Value = [double]"7.869792681105293"
$ws.Range("D3").Value = [double]"-10.45921331889471"
$ws.Range("E3").Value = [double]"-2.319131318894706"
$ws.Range("F3").Value = [double]"-0.5308223188947059"
$ws.Range("G3").Value = [double]"-3.708941318894706"
$ws.Range("H3").Value = [double]"-4.004270318894706"
$ws.Range("I3").Value = [double]"-1.489825318894706"
$ws.Range("J3").Value = [double]"-2.419547318894706"
$ws.Range("K3").Value = [double]"-2.184539318894706"

# Row 4: insert new first-quarter error in B4; shift prior values right by one column
$ws.Range("B4").Value = [double]"0.4006423502480008"
$ws.Range("C4").Value = [double]"-17.928363649752"
$ws.Range("D4").Value = [double]"-9.788281649751999"
$ws.Range("E4").Value = [double]"-7.999972649751999"
$ws.Range("F4").Value = [double]"-11.178091649752"
$ws.Range("G4").Value = [double]"-11.473420649752"
$ws.Range("H4").Value = [double]"-8.958975649751999"
$ws.Range("I4").Value = [double]"-9.888697649751998"
$ws.Range("J4").Value = [double]"-9.653689649752"
$ws.Range("K4").Value = [double]"-9.960825649752"

# Row 5: insert new first-quarter error in B5; shift prior values right by one column
$ws.Range("B5").Value = [double]"0.2428981503854697"
$ws.Range("C5").Value = [double]"8.382980150385469"
$ws.Range("D5").Value = [double]"10.17128915038547"
$ws.Range("E5").Value = [double]"6.99317015038547"
$ws.Range("F5").Value = [double]"6.69784115038547"
$ws.Range("G5").Value = [double]"9.21228615038547"
$ws.Range("H5").Value = [double]"8.28256415038547"
$ws.Range("I5").Value = [double]"8.517572150385469"
$ws.Range("J5").Value = [double]"8.210436150385469"
$ws.Range("K5").Value = [double]"8.721112150385469"

# Row 6: insert new first-quarter error in B6; shift prior values right by one column
$ws.Range("B6").Value = [double]"0.2433086034063205"
$ws.Range("C6").Value = [double]"2.03161760340632"
$ws.Range("D6").Value = [double]"-1.146501396593679"
$ws.Range("E6").Value = [double]"-1.441830396593679"
$ws.Range("F6").Value = [double]"1.072614603406321"
$ws.Range("G6").Value = [double]"0.1428926034063205"
$ws.Range("H6").Value = [double]"0.3779006034063205"
$ws.Range("I6").Value = [double]"0.07076460340632051"
$ws.Range("J6").Value = [double]"0.5814406034063205"
$ws.Range("K6").Value = [double]"0.2940736034063205"

# Row 7: insert new first-quarter error in B7; shift prior values right by one column
$ws.Range("B7").Value = [double]"-0.115952868393072"
$ws.Range("C7").Value = [double]"-3.294071868393072"
$ws.Range("D7").Value = [double]"-3.589400868393072"
$ws.Range("E7").Value = [double]"-1.074955868393072"
$ws.Range("F7").Value = [double]"-2.004677868393072"
$ws.Range("G7").Value = [double]"-1.769669868393072"
$ws.Range("H7").Value = [double]"-2.076805868393072"
$ws.Range("I7").Value = [double]"-1.566129868393072"
$ws.Range("J7").Value = [double]"-1.853496868393072"
$ws.Range("K7").Value = [double]"-1.825497868393072"

# Row 8: insert new first-quarter error in B8; shift prior values right by one column
$ws.Range("B8").Value = [double]"0.1459581181302581"
$ws.Range("C8").Value = [double]"-0.1493708818697419"
$ws.Range("D8").Value = [double]"2.365074118130258"
$ws.Range("E8").Value = [double]"1.435352118130258"
$ws.Range("F8").Value = [double]"1.670360118130258"
$ws.Range("G8").Value = [double]"1.363224118130258"
$ws.Range("H8").Value = [double]"1.873900118130258"
$ws.Range("I8").Value = [double]"1.586533118130258"
$ws.Range("J8").Value = [double]"1.614532118130258"
$ws.Range("K8").Value = [double]"1.707180118130258"

# Row 9: insert new first-quarter error in B9; shift prior values right by one column
$ws.Range("B9").Value = [double]"-0.08819670345554087"
$ws.Range("C9").Value = [double]"2.426248296544459"
$ws.Range("D9").Value = [double]"1.496526296544459"
$ws.Range("E9").Value = [double]"1.731534296544459"
$ws.Range("F9").Value = [double]"1.424398296544459"
$ws.Range("G9").Value = [double]"1.935074296544459"
$ws.Range("H9").Value = [double]"1.647707296544459"
$ws.Range("I9").Value = [double]"1.675706296544459"
$ws.Range("J9").Value = [double]"1.768354296544459"
$ws.Range("K9").Value = [double]"1.981330296544459"

# Row 10: insert new first-quarter error in B10; shift prior values right by one column
$ws.Range("B10").Value = [double]"0.3847923593882046"
$ws.Range("C10").Value = [double]"-0.5449296406117954"
$ws.Range("D10").Value = [double]"-0.3099216406117954"
$ws.Range("E10").Value = [double]"-0.6170576406117954"
$ws.Range("F10").Value = [double]"-0.1063816406117954"
$ws.Range("G10").Value = [double]"-0.3937486406117954"
$ws.Range("H10").Value = [double]"-0.3657496406117954"
$ws.Range("I10").Value = [double]"-0.2731016406117954"
$ws.Range("J10").Value = [double]"-0.06012564061179543"
$ws.Range("K10").Value = [double]"-0.5607856406117955"

# Row 11: insert new first-quarter error in B11; shift prior values right by one column
$ws.Range("B11").Value = [double]"0.03766489642184559"
$ws.Range("C11").Value = [double]"0.2726728964218456"
$ws.Range("D11").Value = [double]"-0.03446310357815441"
$ws.Range("E11").Value = [double]"0.4762128964218456"
$ws.Range("F11").Value = [double]"0.1888458964218456"
$ws.Range("G11").Value = [double]"0.2168448964218456"
$ws.Range("H11").Value = [double]"0.3094928964218456"
$ws.Range("I11").Value = [double]"0.5224688964218456"
$ws.Range("J11").Value = [double]"0.02180889642184558"
$ws.Range("K11").Value = [double]"0.3115658964218456"

# Row 12: insert new first-quarter error in B12; shift prior values right by one column
$ws.Range("B12").Value = [double]"0.1769978556124878"
$ws.Range("C12").Value = [double]"-0.1301381443875122"
$ws.Range("D12").Value = [double]"0.3805378556124878"
$ws.Range("E12").Value = [double]"0.09317085561248779"
$ws.Range("F12").Value = [double]"0.1211698556124878"
$ws.Range("G12").Value = [double]"0.2138178556124878"
$ws.Range("H12").Value = [double]"0.4267938556124878"
$ws.Range("I12").Value = [double]"-0.07386614438751221"
$ws.Range("J12").Value = [double]"0.2158908556124878"
$ws.Range("K12").Value = [double]"-0.04200114438751221"

# Row 13: insert new first-quarter error in B13; shift prior values right by one column
$ws.Range("B13").Value = [double]"0.1302808926112106"
$ws.Range("C13").Value = [double]"0.6409568926112106"
$ws.Range("D13").Value = [double]"0.3535898926112106"
$ws.Range("E13").Value = [double]"0.3815888926112106"
$ws.Range("F13").Value = [double]"0.4742368926112106"
$ws.Range("G13").Value = [double]"0.6872128926112107"
$ws.Range("H13").Value = [double]"0.1865528926112106"
$ws.Range("I13").Value = [double]"0.4763098926112106"
$ws.Range("J13").Value = [double]"0.2184178926112106"
$ws.Range("K13").Value = [double]"0.6007988926112107"

# Row 14: insert new first-quarter error in B14; shift prior values right by one column
$ws.Range("B14").Value = [double]"-0.1944981035472806"
$ws.Range("C14").Value = [double]"-0.4818651035472806"
$ws.Range("D14").Value = [double]"-0.4538661035472806"
$ws.Range("E14").Value = [double]"-0.3612181035472806"
$ws.Range("F14").Value = [double]"-0.1482421035472806"
$ws.Range("G14").Value = [double]"-0.6489021035472806"
$ws.Range("H14").Value = [double]"-0.3591451035472806"
$ws.Range("I14").Value = [double]"-0.6170371035472806"
$ws.Range("J14").Value = [double]"-0.2346561035472806"
$ws.Range("K14").Value = [double]"-0.6459331975472806"

# Row 15: insert new first-quarter error in B15; shift prior values right by one column
$ws.Range("B15").Value = [double]"-0.3817857436446591"
$ws.Range("C15").Value = [double]"-0.3537867436446591"
$ws.Range("D15").Value = [double]"-0.2611387436446591"
$ws.Range("E15").Value = [double]"-0.04816274364465911"
$ws.Range("F15").Value = [double]"-0.5488227436446591"
$ws.Range("G15").Value = [double]"-0.2590657436446591"
$ws.Range("H15").Value = [double]"-0.5169577436446591"
$ws.Range("I15").Value = [double]"-0.1345767436446591"
$ws.Range("J15").Value = [double]"-0.5458538376446591"
$ws.Range("K15").Value = [double]"-0.2588977436446591"

# Row 16: insert new first-quarter error in B16; shift prior values right by one column
$ws.Range("B16").Value = [double]"-2.375649628613696E-07"
$ws.Range("C16").Value = [double]"0.09264776243503714"
$ws.Range("D16").Value = [double]"0.3056237624350371"
$ws.Range("E16").Value = [double]"-0.1950362375649629"
$ws.Range("F16").Value = [double]"0.09472076243503715"
$ws.Range("G16").Value = [double]"-0.1631712375649629"
$ws.Range("H16").Value = [double]"0.2192097624350371"
$ws.Range("I16").Value = [double]"-0.1920673315649629"
$ws.Range("J16").Value = [double]"0.09488876243503713"

# Row 17: insert new first-quarter error in B17; shift prior values right by one column
$ws.Range("B17").Value = [double]"-0.0555296279974082"
$ws.Range("C17").Value = [double]"0.1574463720025918"
$ws.Range("D17").Value = [double]"-0.3432136279974082"
$ws.Range("E17").Value = [double]"-0.05345662799740819"
$ws.Range("F17").Value = [double]"-0.3113486279974082"
$ws.Range("G17").Value = [double]"0.0710323720025918"
$ws.Range("H17").Value = [double]"-0.3402447219974082"
$ws.Range("I17").Value = [double]"-0.0532886279974082"

# Row 18: insert new first-quarter error in B18; shift prior values right by one column
$ws.Range("B18").Value = [double]"3.829984367986761E-07"
$ws.Range("C18").Value = [double]"-0.5006596170015631"
$ws.Range("D18").Value = [double]"-0.2109026170015632"
$ws.Range("E18").Value = [double]"-0.4687946170015632"
$ws.Range("F18").Value = [double]"-0.08641361700156319"
$ws.Range("G18").Value = [double]"-0.4976907110015631"
$ws.Range("H18").Value = [double]"-0.2107346170015632"

# Row 19: insert new first-quarter error in B19; shift prior values right by one column
$ws.Range("B19").Value = [double]"-1.604754923945073E-07"
$ws.Range("C19").Value = [double]"0.2897568395245076"
$ws.Range("D19").Value = [double]"0.03186483952450761"
$ws.Range("E19").Value = [double]"0.4142458395245076"
$ws.Range("F19").Value = [double]"0.002968745524507627"
$ws.Range("G19").Value = [double]"0.2899248395245076"

# Row 20: insert new first-quarter error in B20; shift prior values right by one column
$ws.Range("B20").Value = [double]"0.009398958989038461"
$ws.Range("C20").Value = [double]"-0.2484930410109615"
$ws.Range("D20").Value = [double]"0.1338879589890384"
$ws.Range("E20").Value = [double]"-0.2773891350109615"
$ws.Range("F20").Value = [double]"0.009566958989038449"

# Row 21: insert new first-quarter error in B21; shift prior values right by one column
$ws.Range("B21").Value = [double]"-0.07651818316594991"
$ws.Range("C21").Value = [double]"0.3058628168340501"
$ws.Range("D21").Value = [double]"-0.1054142771659499"
$ws.Range("E21").Value = [double]"0.1815418168340501"

# Row 22: insert new first-quarter error in B22; shift prior values right by one column
$ws.Range("B22").Value = [double]"2.770877186031306E-07"
$ws.Range("C22").Value = [double]"-0.4112768169122814"
$ws.Range("D22").Value = [double]"-0.1243207229122814"

# Row 23: insert new first-quarter error in B23; shift prior values right by one column
$ws.Range("B23").Value = [double]"0.2010531357750048"
$ws.Range("C23").Value = [double]"0.4880092297750048"

# Row 24: insert new first-quarter error in B24; shift prior values right by one column
$ws.Range("B24").Value = [double]"-0.2003621554241067"

